$wb = $excel.ActiveWorkbook

# Hunk 0: sheet ALC, row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2190
$ws.Range("J40").Value = 2190
$ws.Range("L40").Value = 2190
$ws.Range("N40").Value = -2540

# Hunk 1: sheet ALC, row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 266.85715
$ws.Range("J41").Value = 338.33334
$ws.Range("L41").Value = 338.33334
$ws.Range("N41").Value = -1218.33334

# Hunk 2: sheet ALC, row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1532.5
$ws.Range("I70").Value = 1373.75
$ws.Range("K70").Value = 4121.25
$ws.Range("M70").Value = -3851.25

# Hunk 3: sheet ALC, row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1532.5
$ws.Range("I73").Value = 1373.75
$ws.Range("K73").Value = 4121.25
$ws.Range("M73").Value = -3185.25

# Hunk 4: sheet ALC, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1905.5714
$ws.Range("I132").Value = 1798.8182
$ws.Range("K132").Value = 5396.4546
$ws.Range("M132").Value = -2866.4546

# Hunk 5: sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3890.3845
$ws.Range("I138").Value = 1476.2941
$ws.Range("J138").Value = 5755.8184
$ws.Range("K138").Value = 4428.8823
$ws.Range("L138").Value = 17267.4552
$ws.Range("M138").Value = 711.1176999999998
$ws.Range("N138").Value = -27547.4552

# Hunk 6: sheet ALC, row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2318.762
$ws.Range("I141").Value = 1931.2106
$ws.Range("K141").Value = 5793.6318
$ws.Range("M141").Value = -613.6318000000001

# Hunk 7: sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1136.8235
$ws.Range("I74").Value = 1132.9375
$ws.Range("K74").Value = 1132.9375
$ws.Range("M74").Value = -258.9375

# Hunk 8: sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1136.8235
$ws.Range("I77").Value = 1132.9375
$ws.Range("K77").Value = 5664.6875
$ws.Range("M77").Value = -1296.6875

# Hunk 9: sheet ARM, row 94
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

# Hunk 10: sheet ARM, row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1500
$ws.Range("I110").Value = 1500
$ws.Range("K110").Value = 1500
$ws.Range("M110").Value = 545

# Hunk 11: sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2070.077
$ws.Range("I122").Value = 2177.0908
$ws.Range("K122").Value = 6531.2724
$ws.Range("M122").Value = -4081.2724

# Hunk 12: sheet BSM, row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

# Hunk 13: sheet BSM, row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# Hunk 14: sheet BSM, row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1421.4
$ws.Range("I99").Value = 1636
$ws.Range("J99").Value = 1099.5
$ws.Range("K99").Value = 1636
$ws.Range("L99").Value = 1099.5
$ws.Range("M99").Value = -138
$ws.Range("N99").Value = -4095.5

# Hunk 15: sheet BSM, row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3222.7778
$ws.Range("I105").Value = 3250.625
$ws.Range("K105").Value = 3250.625
$ws.Range("M105").Value = -1503.625

# Hunk 16: sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4083.4
$ws.Range("I134").Value = 4231.5557
$ws.Range("J134").Value = 2750
$ws.Range("K134").Value = 12694.6671
$ws.Range("L134").Value = 8250
$ws.Range("M134").Value = -10159.6671
$ws.Range("N134").Value = -13320

# Hunk 17: sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2497.6
$ws.Range("I31").Value = 2926.2856
$ws.Range("J31").Value = 1497.3334
$ws.Range("K31").Value = 2926.2856
$ws.Range("L31").Value = 1497.3334
$ws.Range("M31").Value = -2631.2856
$ws.Range("N31").Value = -2087.3334

# Hunk 18: sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2497.6
$ws.Range("I34").Value = 2926.2856
$ws.Range("J34").Value = 1497.3334
$ws.Range("K34").Value = 2926.2856
$ws.Range("L34").Value = 1497.3334
$ws.Range("M34").Value = -2724.2856
$ws.Range("N34").Value = -1901.3334

# Hunk 19: sheet CRP, row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2239.2
$ws.Range("I122").Value = 2239.2
$ws.Range("K122").Value = 6717.599999999999
$ws.Range("M122").Value = -4267.599999999999

# Hunk 20: sheet CUL, row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 137528780
$ws.Range("I4").Value = 104880180
$ws.Range("J4").Value = 170177400
$ws.Range("K4").Value = 314640540
$ws.Range("L4").Value = 510532200
$ws.Range("M4").Value = -314640428
$ws.Range("N4").Value = -510532424

# Hunk 21: sheet CUL, row 9
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 10752
$ws.Range("J9").Value = 10752
$ws.Range("L9").Value = 32256
$ws.Range("N9").Value = -32704

# Hunk 22: sheet CUL, row 10
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 100
$ws.Range("I10").Value = 100
$ws.Range("K10").Value = 300
$ws.Range("M10").Value = -161

# Hunk 23: sheet CUL, row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 244.66667
$ws.Range("I12").Value = 235.77777
$ws.Range("J12").Value = 258
$ws.Range("K12").Value = 707.33331
$ws.Range("L12").Value = 774
$ws.Range("M12").Value = -534.33331
$ws.Range("N12").Value = -1120

# Hunk 24: sheet CUL, row 13
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 1283
$ws.Range("I13").Value = 1283
$ws.Range("K13").Value = 3849
$ws.Range("M13").Value = -3681

# Hunk 25: sheet CUL, row 15
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()

# Hunk 26: sheet CUL, row 23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 324.5
$ws.Range("I23").Value = 300
$ws.Range("K23").Value = 900
$ws.Range("M23").Value = -665

# Hunk 27: sheet CUL, row 86
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 293.25
$ws.Range("J86").Value = 289
$ws.Range("L86").Value = 867
$ws.Range("N86").Value = -3239

# Hunk 28: sheet CUL, row 89
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 293.25
$ws.Range("J89").Value = 289
$ws.Range("L89").Value = 2601
$ws.Range("N89").Value = -14457

# Hunk 29: sheet CUL, row 124
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 2030
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

# Hunk 30: sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1556.8077
$ws.Range("I131").Value = 949.8570999999999
$ws.Range("J131").Value = 1780.421
$ws.Range("K131").Value = 2849.5713
$ws.Range("L131").Value = 5341.263
$ws.Range("M131").Value = 2190.4287
$ws.Range("N131").Value = -15421.263

# Hunk 31: sheet CUL, row 141
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 14154.833
$ws.Range("I141").Value = 3309.6667
$ws.Range("K141").Value = 9929.000100000001
$ws.Range("M141").Value = -4749.000100000001

# Hunk 32: sheet GSM, row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11429.5
$ws.Range("I70").Value = 10770.857
$ws.Range("K70").Value = 10770.857
$ws.Range("M70").Value = -10500.857

# Hunk 33: sheet GSM, row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 11429.5
$ws.Range("I73").Value = 10770.857
$ws.Range("K73").Value = 10770.857
$ws.Range("M73").Value = -9834.857

# Hunk 34: sheet GSM, row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5838.2
$ws.Range("I102").Value = 5798.3335
$ws.Range("J102").Value = 5898
$ws.Range("K102").Value = 5798.3335
$ws.Range("L102").Value = 5898
$ws.Range("M102").Value = -4176.3335
$ws.Range("N102").Value = -9142

# Hunk 35: sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1475
$ws.Range("I132").Value = 1475
$ws.Range("K132").Value = 4425
$ws.Range("M132").Value = -1895

# Hunk 36: sheet LTW, row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 940
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

# Hunk 37: sheet LTW, row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 940
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

# Hunk 38: sheet LTW, row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3720.4
$ws.Range("I40").Value = 3400.75
$ws.Range("K40").Value = 3400.75
$ws.Range("M40").Value = -3264.75

# Hunk 39: sheet LTW, row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1874.25
$ws.Range("I82").Value = 1874.25
$ws.Range("K82").Value = 1874.25
$ws.Range("M82").Value = -1513.25

# Hunk 40: sheet LTW, row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1874.25
$ws.Range("I85").Value = 1874.25
$ws.Range("K85").Value = 1874.25
$ws.Range("M85").Value = -626.25

# Hunk 41: sheet LTW, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2955.7778
$ws.Range("I136").Value = 2869.5652
$ws.Range("J136").Value = 3451.5
$ws.Range("K136").Value = 8608.695599999999
$ws.Range("L136").Value = 10354.5
$ws.Range("M136").Value = -6058.695599999999
$ws.Range("N136").Value = -15454.5

# Hunk 42: sheet LTW, row 141
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H141").Value = 116633.336
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 116633.336
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 116633.336
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -126993.336

# Hunk 43: sheet WVR, row 3
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 5000
$ws.Range("J3").Value = 5000
$ws.Range("L3").Value = 5000
$ws.Range("N3").Value = -5228

# Hunk 44: sheet WVR, row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31248

# Hunk 45: sheet WVR, row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156240

# Hunk 46: sheet WVR, row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9874.111000000001
$ws.Range("I81").Value = 1372.6
$ws.Range("J81").Value = 20501
$ws.Range("K81").Value = 2745.2
$ws.Range("L81").Value = 41002
$ws.Range("M81").Value = -1684.2
$ws.Range("N81").Value = -43124

# Hunk 47: sheet WVR, row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 9874.111000000001
$ws.Range("I84").Value = 1372.6
$ws.Range("J84").Value = 20501
$ws.Range("K84").Value = 13726
$ws.Range("L84").Value = 205010
$ws.Range("M84").Value = -8422
$ws.Range("N84").Value = -215618

# Hunk 48: sheet WVR, row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 404.25
$ws.Range("I107").Value = 348
$ws.Range("K107").Value = 1044
$ws.Range("M107").Value = 876

# Hunk 49: sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1854.3636
$ws.Range("I122").Value = 1933
$ws.Range("J122").Value = 1500.5
$ws.Range("K122").Value = 5799
$ws.Range("L122").Value = 4501.5
$ws.Range("M122").Value = -3349
$ws.Range("N122").Value = -9401.5

# Hunk 50: sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2296.3845
$ws.Range("I132").Value = 1742.6111
$ws.Range("K132").Value = 5227.8333
$ws.Range("M132").Value = -2697.8333

# Hunk 51: sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1503.875
$ws.Range("J136").Value = 2342.3333
$ws.Range("L136").Value = 7026.999899999999
$ws.Range("N136").Value = -12126.9999
